$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the sheet
$ws.Name = "ValidLogin"

# Add the data
$ws.Range("A1").Value = "UserName"
$ws.Range("B1").Value = "Password"
$ws.Range("A2").Value = "admin"
$ws.Range("B2").Value = "manager"

# Set the selection / active cell as per the diff
$ws.Range("C2").Select()

# Zoom settings
$ws.Application.ActiveWindow.Zoom = 235
